$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AA, shifting the existing "nom" (AA) and
# "url_produit" (AB) columns one to the right (to AB and AC respectively).
$ws.Columns("AA:AA").Insert()

# Header cell for the freshly inserted column: a new timestamp column,
# matching the pattern of the other date/time header cells.
$ws.Range("AA1").Value = "2026-01-28 20:16:50"

# Populate the new column for each data row (2-205) with the latest
# observed value, mirroring column Z (the previous last price-history
# column) for that row. Rows without a price in Z stay blank, but we still
# materialize the cell (matching Z's own blank-but-present cell) by
# copying Z's formatting across.
for ($r = 2; $r -le 205; $r++) {
    $zVal = $ws.Cells.Item($r, 26).Value()
    if ($zVal -ne $null -and $zVal -ne "") {
        $ws.Cells.Item($r, 27).Value = $zVal
    } else {
        $ws.Cells.Item($r, 26).Copy()
        $ws.Cells.Item($r, 27).PasteSpecial(-4122)
    }
}
